$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 187, pushing the
# previously existing rows 187-202 down to 188-203 (dimension A1:R202 -> A1:R203).
$ws.Rows("187:187").Insert()

$ws.Range("A187").Value = 3
$ws.Range("B187").Value = 'Femacal de La Calera'
$ws.Range("C187").Value = 'Coquimbo'
$ws.Range("D187").Value = 44826
$ws.Range("E187").Value = 5
$ws.Range("F187").Value = 100112026
$ws.Range("G187").Value = 'Haba'
$ws.Range("H187").Value = 'Sin especificar'
$ws.Range("I187").Value = 'Primera'
$ws.Range("J187").Value = 105
$ws.Range("K187").Value = 12000
$ws.Range("L187").Value = 12500
$ws.Range("M187").Value = 12262
$ws.Range("N187").Value = '$/malla 25 kilos'
$ws.Range("O187").Value = 'Provincia de Limarí'
$ws.Range("P187").Value = 490
$ws.Range("Q187").Value = 25
$ws.Range("R187").Value = 'Hortaliza'
